# Adding QT TA data: two new rows (148-149) of CRM196 titration data to the
# CRMAccuracyData sheet, including a note about a failed titration run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 148 ---------------------------------------------------------
# Fill the text/label columns first (F, then G, then B/D) so that the new
# shared-string entries land in the same order they were originally typed:
#   50 -> CRM196_opened20240828
#   51 -> Issue with titration, complete values were not recorded in file
#   52 -> N/A
$ws.Range("F148").Value2 = "CRM196_opened20240828"
$ws.Range("G148").Value2 = "Issue with titration, complete values were not recorded in file"
$ws.Range("B148").Value2 = "N/A"
$ws.Range("D148").Value2 = "N/A"

# Numeric columns for row 148
$ws.Range("A148").Value2 = 20240924
$ws.Range("C148").Value2 = 2215.3200000000002
$ws.Range("E148").Value2 = 196

# --- Row 149 ---------------------------------------------------------
$ws.Range("A149").Value2 = 20240924
$ws.Range("C149").Value2 = 2215.3200000000002
$ws.Range("E149").Value2 = 196
$ws.Range("F149").Value2 = "CRM196_opened20240829"

# --- Formatting: match the existing data rows (style index used by the
# rest of the table, i.e. 12pt font with no explicit color override) by
# copying formats from the last existing data row (147) for the columns
# that carry that style (A, C, E, F). Columns B, D, G stay unstyled, same
# as in row 148 for this entry.
$ws.Range("A147").Copy() | Out-Null
$ws.Range("A148").PasteSpecial(-4122) | Out-Null
$ws.Range("C147").Copy() | Out-Null
$ws.Range("C148").PasteSpecial(-4122) | Out-Null
$ws.Range("E147").Copy() | Out-Null
$ws.Range("E148").PasteSpecial(-4122) | Out-Null
$ws.Range("F147").Copy() | Out-Null
$ws.Range("F148").PasteSpecial(-4122) | Out-Null

$ws.Range("A147").Copy() | Out-Null
$ws.Range("A149").PasteSpecial(-4122) | Out-Null
$ws.Range("C147").Copy() | Out-Null
$ws.Range("C149").PasteSpecial(-4122) | Out-Null
$ws.Range("E147").Copy() | Out-Null
$ws.Range("E149").PasteSpecial(-4122) | Out-Null
$ws.Range("F147").Copy() | Out-Null
$ws.Range("F149").PasteSpecial(-4122) | Out-Null

# --- View state: leave the cursor on the next empty entry row, scrolled
# down near the bottom of the table, similar to where the author left off.
$win = $excel.ActiveWindow
$win.ScrollRow = 130
$win.ScrollColumn = 1
$ws.Range("B149").Select() | Out-Null
